$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: row numbers for the new rows (4-12) ---
for ($r = 4; $r -le 12; $r++) {
    $ws.Range("A$r").Value = $r
}

# --- Column C: Bloom's taxonomy level codes for rows 4-12, filled in row order ---
# (this introduces new shared strings K2, K4, K1 in this exact first-seen order)
$ws.Range("C4").Value = "K2"
$ws.Range("C5").Value = "K3"
$ws.Range("C6").Value = "K4"
$ws.Range("C7").Value = "K3"
$ws.Range("C8").Value = "K3"
$ws.Range("C9").Value = "K3"
$ws.Range("C10").Value = "K2"
$ws.Range("C11").Value = "K1"
$ws.Range("C12").Value = "K2"

# --- Column B: question text ---
# Row 1 keeps its question, mark reduced from (15) to (8)
$ws.Range("B1").Value = "User ‘A’ in delhi wishes to send a file for printing to user ‘B’ in florida, whose system is connected to a printer; while user ‘C’ from tokyo wants to save a video file in the hard disk of user ‘D’ in london. Model the required peer-to-peer network architecture.(8)"
# Row 3 ("Develop...") is (re)written before row 2 ("Explain a formal proof...")
# so the shared-string table lists them in that order, matching the source file.
$ws.Range("B3").Value = "Develop a detailed implementation of causal consistency, and provide a correctness argument for your implementation.(8)"
$ws.Range("B2").Value = "Explain a formal proof to justify the correctness of algorithm that construct sequential consistency using local read operations.(8)"
$ws.Range("B4").Value = "Explain the concept of distributed systems and discuss the advantages and challenges associated with them.(8)"
$ws.Range("B5").Value = "What is the role of middleware in a distributed system? Describe some common middleware technologies used in distributed systems.(8)"
$ws.Range("B6").Value = "Discuss the client-server architecture in the context of distributed systems. Explain the responsibilities of the client and server components.(8)"
$ws.Range("B7").Value = "What is the difference between synchronous and asynchronous communication in distributed systems? Provide examples of each.(8)"
$ws.Range("B8").Value = "Describe the different types of distributed system architectures, such as peer-to-peer, client-server, and hybrid architectures. Compare and contrast their characteristics.(8)"
$ws.Range("B9").Value = "Explain the concept of fault tolerance in distributed systems. Discuss various techniques used to achieve fault tolerance, such as replication and redundancy.(8)"
$ws.Range("B10").Value = "What is distributed file system? Discuss the design principles and features of a distributed file system.(8)"
$ws.Range("B11").Value = "Describe the challenges of data consistency in distributed systems. Explain the concepts of eventual consistency and strong consistency.(8)"
$ws.Range("B12").Value = "Discuss the role of distributed transactions in ensuring data integrity across multiple distributed components. Explain the ACID properties and their significance in distributed transactions.(8)"

# --- Column D: course code for the new rows (4-12) ---
for ($r = 4; $r -le 12; $r++) {
    $ws.Range("D$r").Value = "C314.5"
}

# --- Column widths ---
# Target stored width is 171.7109375 characters; the engine (like real Excel)
# snaps ColumnWidth to its internal pixel grid, so we feed it the input that
# quantizes closest to the target (171.7109375 itself would round further away).
$ws.Columns.Item(2).ColumnWidth = 170.8

# --- View: zoom level and active selection ---
$excel.ActiveWindow.Zoom = 90
$ws.Range("B17").Select()
